$wb = $excel.ActiveWorkbook

# Rename sheets to reflect new "01_" numbering scheme
$wb.Worksheets.Item("02_1_general").Name = "01_1_general"
$wb.Worksheets.Item("02_1_diccionario").Name = "01_1_diccionario"

$ws = $wb.Worksheets.Item("01_1_diccionario")

# Insert a new "Nombre propuesto" column right after "Nombre de la variable" (col A)
$ws.Range("B1").EntireColumn.Insert()

# The old "Unidad de Medida" column (now shifted to column E) is no longer used
$ws.Range("E1").EntireColumn.Delete()

# Update header row text
$ws.Range("A1").Value = "Nombre actual"
$ws.Range("B1").Value = "Nombre propuesto"

# Populate the new column with a proposed (lower-case) variable name
$ws.Range("B2").Formula = "=LOWER(A2)"
$ws.Range("B3:B18").Formula = "=LOWER(A3)"

$ws.Columns.AutoFit()

$ws.Range("B23").Select()
